$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bareme")

# Update the data cells (rows 3-32) on sheet "Bareme".
# Row 2 contains SUM() formulas over each column (rows 3:145) and will
# recalculate automatically once the underlying values change.

$ws.Range("M3").Value = 2
$ws.Range("Z3").Value = 2
$ws.Range("N4").Value = 3
$ws.Range("AA4").Value = 3
$ws.Range("O5").Value = 1
$ws.Range("AB5").Value = 1
$ws.Range("P6").Value = 2
$ws.Range("AC6").Value = 2
$ws.Range("Q7").Value = 4
$ws.Range("V7").Value = 3
$ws.Range("AD7").Value = 4
$ws.Range("W8").Value = 1
$ws.Range("X9").Value = 2
$ws.Range("T10").Value = 1
$ws.Range("Y10").Value = 4
$ws.Range("AG10").Value = 1
$ws.Range("U11").Value = 2
$ws.Range("M12").Value = 0
$ws.Range("AA13").Value = 3
$ws.Range("D15").Value = 1
$ws.Range("M15").Value = 2
$ws.Range("P15").Value = 0
$ws.Range("Z15").Value = 2
$ws.Range("E16").Value = 2
$ws.Range("N16").Value = 3
$ws.Range("Q16").Value = 0
$ws.Range("AA16").Value = 3
$ws.Range("F17").Value = 3
$ws.Range("O17").Value = 1
$ws.Range("R17").Value = 0
$ws.Range("G18").Value = 1
$ws.Range("P18").Value = 2
$ws.Range("S18").Value = 0
$ws.Range("H19").Value = 2
$ws.Range("Q19").Value = 4
$ws.Range("T19").Value = 0
$ws.Range("V19").Value = 3
$ws.Range("I20").Value = 4
$ws.Range("U20").Value = 0
$ws.Range("W20").Value = 1
$ws.Range("V21").Value = 0
$ws.Range("X21").Value = 2
$ws.Range("T22").Value = 1
$ws.Range("W22").Value = 0
$ws.Range("Y22").Value = 4
$ws.Range("L23").Value = 1
$ws.Range("U23").Value = 2
$ws.Range("X23").Value = 0
$ws.Range("Y24").Value = 0
$ws.Range("N25").Value = 3
$ws.Range("Z25").Value = 0
$ws.Range("AB25").Value = 1
$ws.Range("AA26").Value = 0
$ws.Range("AC26").Value = 2
$ws.Range("P27").Value = 1
$ws.Range("AB27").Value = 0
$ws.Range("Q28").Value = 2
$ws.Range("V28").Value = 3
$ws.Range("AC28").Value = 0
$ws.Range("R29").Value = 3
$ws.Range("AD29").Value = 0
$ws.Range("S30").Value = 1
$ws.Range("AE30").Value = 0
$ws.Range("T31").Value = 1
$ws.Range("AF31").Value = 0
$ws.Range("U32").Value = 1
$ws.Range("AG32").Value = 0

# Force a full recalculation so the SUM formulas in row 2 pick up the new values.
$excel.CalculateFullRebuild()

# Update the view state: scroll so column B is the left-most visible column,
# and move the active selection to AB5.
$ws.Activate()
$ws.Range("AB5").Select()
$excel.ActiveWindow.ScrollColumn = 2
